# Finishes conversion of nodes to new format.
#
# The "other" worksheet tracked a set of misc. Clava AST node kinds. This
# edit removes the rows for node kinds that no longer apply / were merged
# during the restructuring (several Comment subtypes, a few stale
# placeholder entries, and some superseded Omp* entries), and marks the
# surviving rows as completed ("o") where appropriate. Removing these rows
# also makes the corresponding shared-string entries unused, so they drop
# out of xl/sharedStrings.xml automatically when the workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("other")

# Rows (1-based, as they exist in the *current* sheet) to remove, deleted
# from the bottom up so earlier indices in the list stay valid as we go:
#   26 OMPParallelDirective
#   24 OMPExecutableDirective
#   23 OmpDirectiveKind
#   21 Undefined
#   16 OriginalNamespace
#   15 NullNodeOld
#   12 TextComment
#   11 ParagraphComment
#    9 InlineContentComment
#    7 InlineCommandComment
#    6 FullComment
#    3 BlockContentComment
$rowsToDelete = @(26, 24, 23, 21, 16, 15, 12, 11, 9, 7, 6, 3)
foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete()
}

# After the deletions above, the remaining "Comment" family rows land at
# 3-6 and need their Status column filled in (Comment itself stays
# half-way "-", the rest are now fully converted "o").
$ws.Cells.Item(3, 2).Value = "-"   # Comment
$ws.Cells.Item(4, 2).Value = "o"   # DummyComment
$ws.Cells.Item(5, 2).Value = "o"   # InlineComment
$ws.Cells.Item(6, 2).Value = "o"   # MultiLineComment

# Reset the view: scroll back to the top and select B7 (was B31/topLeftCell
# A7 while the sheet had more rows).
$ws.Range("B7").Select() | Out-Null
